$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "22" + " Apr 2020: " (two runs) -> single run "22 Apr 2020: "
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("22 Apr 2020: ", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "22 Apr 2020: ", 2)

# ---------------------------------------------------------------------------
# Change 2: append the 30 Apr 2020 notes before the trailing bookmark
# paragraph (the very last, empty paragraph that holds the _GoBack bookmark).
# ---------------------------------------------------------------------------
$anchor = $d.Paragraphs.Last.Range

$lines = @(
    "",
    "",
    "30 Apr 2020:",
    "",
    "PLXData object successfully reads in data from .plx files and stores in a struct, P.",
    "Data organization is completely different from the exported mat  format file from OFS.",
    [string]::Format("Two options: work with data in P, or {0}reexport{1} data in P to the mat format.", [char]0x201C, [char]0x201D),
    "",
    "Advantages: mat format already has some code to pull out channels, spikes into curves.",
    [string]::Format("Will need to write export/conversion code (not a major deal). But doesn{0}t seem to have PCA info. Necessary? Could either have separate method for spike data, OR have dummy columns", [char]0x2019),
    ""
)

$block = ($lines -join "`r") + "`r"

$anchor.InsertBefore($block)
